$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.06082199999999999
$ws.Range("I2").Value = 0.1716860072883705
$ws.Range("J2").Value = 0.1716860072883705
$ws.Range("M2").Value = 78.99738599999999
$ws.Range("N2").Value = 236.992158
$ws.Range("O2").Value = 0.6882088488047822
$ws.Range("P2").Value = 0.6882088488047821
$ws.Range("Q2").Value = 4.804779011291999
$ws.Range("R2").Value = 43.24301110162799
$ws.Range("S2").Value = 0.1181558294318189
$ws.Range("T2").Value = 0.1181558294318189

# Row 3
$ws.Range("G3").Value = 0.06082199999999999
$ws.Range("I3").Value = 0.1716860072883705
$ws.Range("J3").Value = 0.1716860072883705
$ws.Range("O3").Value = 0.1184214915836591
$ws.Range("P3").Value = 0.118421491583659
$ws.Range("Q3").Value = 0.8267680635539999
$ws.Range("R3").Value = 7.440912571985999
$ws.Range("S3").Value = 0.0203313130671318
$ws.Range("T3").Value = 0.02033131306713179

# Row 4
$ws.Range("G4").Value = 0.06082199999999999
$ws.Range("I4").Value = 0.1716860072883705
$ws.Range("J4").Value = 0.1716860072883705
$ws.Range("M4").Value = 11.10084966666667
$ws.Range("N4").Value = 33.302549
$ws.Range("O4").Value = 0.09670830082721493
$ws.Range("P4").Value = 0.0967083008272149
$ws.Range("Q4").Value = 0.6751758784259999
$ws.Range("R4").Value = 6.076582905833999
$ws.Range("S4").Value = 0.01660346204066715
$ws.Range("T4").Value = 0.01660346204066715

# Row 5
$ws.Range("G5").Value = 0.06082199999999999
$ws.Range("I5").Value = 0.1716860072883705
$ws.Range("J5").Value = 0.1716860072883705
$ws.Range("M5").Value = 2.440340333333334
$ws.Range("N5").Value = 7.321021000000001
$ws.Range("O5").Value = 0.02125973904370978
$ws.Range("P5").Value = 0.02125973904370977
$ws.Range("Q5").Value = 0.148426379754
$ws.Range("R5").Value = 1.335837417786
$ws.Range("S5").Value = 0.003649999712407212
$ws.Range("T5").Value = 0.003649999712407211

# Row 6
$ws.Range("G6").Value = 0.06082199999999999
$ws.Range("I6").Value = 0.1716860072883705
$ws.Range("J6").Value = 0.1716860072883705
$ws.Range("M6").Value = 8.655121
$ws.Range("N6").Value = 25.965363
$ws.Range("O6").Value = 0.07540161974063414
$ws.Range("P6").Value = 0.07540161974063413
$ws.Range("Q6").Value = 0.5264217694619999
$ws.Range("R6").Value = 4.737795925157999
$ws.Range("S6").Value = 0.01294540303634545
$ws.Range("T6").Value = 0.01294540303634545

# Row 7
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.293441
$ws.Range("H7").Value = 0.880323
$ws.Range("I7").Value = 0.8283139927116295
$ws.Range("J7").Value = 0.8283139927116295
$ws.Range("M7").Value = 78.99738599999999
$ws.Range("N7").Value = 236.992158
$ws.Range("O7").Value = 0.6882088488047822
$ws.Range("P7").Value = 0.6882088488047821
$ws.Range("Q7").Value = 23.181071945226
$ws.Range("R7").Value = 208.629647507034
$ws.Range("S7").Value = 0.5700530193729633
$ws.Range("T7").Value = 0.5700530193729632

# Row 8
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.293441
$ws.Range("H8").Value = 0.880323
$ws.Range("I8").Value = 0.8283139927116295
$ws.Range("J8").Value = 0.8283139927116295
$ws.Range("O8").Value = 0.1184214915836591
$ws.Range("P8").Value = 0.118421491583659
$ws.Range("Q8").Value = 3.988814036653666
$ws.Range("R8").Value = 35.899326329883
$ws.Range("S8").Value = 0.09809017851652728
$ws.Range("T8").Value = 0.09809017851652725

# Row 9
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.293441
$ws.Range("H9").Value = 0.880323
$ws.Range("I9").Value = 0.8283139927116295
$ws.Range("J9").Value = 0.8283139927116295
$ws.Range("M9").Value = 11.10084966666667
$ws.Range("N9").Value = 33.302549
$ws.Range("O9").Value = 0.09670830082721493
$ws.Range("P9").Value = 0.0967083008272149
$ws.Range("Q9").Value = 3.257444427036333
$ws.Range("R9").Value = 29.316999843327
$ws.Range("S9").Value = 0.08010483878654778
$ws.Range("T9").Value = 0.08010483878654777

# Row 10
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.293441
$ws.Range("H10").Value = 0.880323
$ws.Range("I10").Value = 0.8283139927116295
$ws.Range("J10").Value = 0.8283139927116295
$ws.Range("M10").Value = 2.440340333333334
$ws.Range("N10").Value = 7.321021000000001
$ws.Range("O10").Value = 0.02125973904370978
$ws.Range("P10").Value = 0.02125973904370977
$ws.Range("Q10").Value = 0.7160959077536668
$ws.Range("R10").Value = 6.444863169783001
$ws.Range("S10").Value = 0.01760973933130257
$ws.Range("T10").Value = 0.01760973933130256

# Row 11
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.293441
$ws.Range("H11").Value = 0.880323
$ws.Range("I11").Value = 0.8283139927116295
$ws.Range("J11").Value = 0.8283139927116295
$ws.Range("M11").Value = 8.655121
$ws.Range("N11").Value = 25.965363
$ws.Range("O11").Value = 0.07540161974063414
$ws.Range("P11").Value = 0.07540161974063413
$ws.Range("Q11").Value = 2.539767361361
$ws.Range("R11").Value = 22.857906252249
$ws.Range("S11").Value = 0.06245621670428869
$ws.Range("T11").Value = 0.06245621670428868
